$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Duplicate the T_AS_S_2 block (A30:E33, with its style/merge/borders) down into
#    rows 38:41 to seed a new 4-row test-case block. Copy(Destination) preserves the
#    per-cell styles (s=4 for A/B/C/E, s=2 for D) and auto-creates the matching
#    vertical merges + extends the sheet dimension, exactly like the source author
#    copy/pasting an existing block as a starting point for a new test case.
$ws.Range("A30:E33").Copy($ws.Range("A38"))

# 2) Turn the old duplicate block at rows 34:37 (currently an exact copy of T_AS_S_2)
#    into the new T_AS_S_3 "Auto.Schedule.Location" test case.
$ws.Range("A34").Value = "T_AS_S_3"
$ws.Range("B34").Value = "Auto.Schedule.Location"
$ws.Range("C34").Value = "Verifies AS schedules events within correct location events"
$ws.Range("D35").Value = "2. Create a VS with location events for the same location as specified in the VD"
$ws.Range("D34").Value = "1. Create a VD with a  location"
$ws.Range("D23").Value = "2. Create an VS"
$ws.Range("E34").Value = "Events with a specified location are contained within their a location event for the correct location."

# 3) Finish filling in the new block (rows 38:41) as the T_AS_S_4 "Auto.Schedule.MoveEvent"
#    test case. (D38/D40/D41 already match the desired text after the copy above.)
$ws.Range("A38").Value = "T_AS_S_3"
$ws.Range("B38").Value = "Auto.Schedule.MoveEvent"
$ws.Range("C38").Value = "Verifies AS can move non-locked, non-location events previously scheduled by the user."
$ws.Range("E38").Value = "Non-locked Events are moved to create space for AS's events so that they don't conflict with other events."
$ws.Range("D39").Value = "2. Create a VS with non-locked and locked events. Events within the VD's valid times must leave only just enough free time for deadline's total work time but not enough with the addition of breaks and max time constraints."

# 4) Row heights to fit the (now longer / shorter) wrapped text in column D.
$ws.Rows(35).RowHeight = 30
$ws.Rows(39).RowHeight = 90
$ws.Rows(40).RowHeight = 30
$ws.Rows(41).RowHeight = 30

# 5) Scroll / selection state left by the editor when they finished authoring the new rows.
$excel.ActiveWindow.ScrollRow = 24
$ws.Range("D40").Select() | Out-Null
